# Generate Report for handoff
#
# The previous handoff attempt (file 13bf0d2b-d7cb-415e-802d-0ddcdb7b57f1.md)
# is replaced by a new attempt (786c73be-f044-4ba5-8f37-ca032f9a9e06.md) whose
# transform failed, so the per-language rows lose their "ready" handoff
# target/date/reason info and report the new status + "Ignored" reason
# instead.

$wb = $excel.ActiveWorkbook

$oldGuidMd  = "13bf0d2b-d7cb-415e-802d-0ddcdb7b57f1.md"
$newGuidMd  = "786c73be-f044-4ba5-8f37-ca032f9a9e06.md"
$oldStatus  = "Ready for handoff"
$newStatus  = "Handoff transform failed"
$mdUrl      = "https://github.com/OpenLocalizationTest/oltest/blob/27f88beec385f9a1e78a100cb74d54f13b0492c1/e2e/$newGuidMd"
$configUrl  = "https://github.com/OpenLocalizationTest/oltest/blob/27f88beec385f9a1e78a100cb74d54f13b0492c1/.localization-config"
$noDate     = "0001-01-01 00:00:00"
$ignored    = "Ignored"
$hlColor    = 15570276   # RGB(100,149,237) == FF6495ED, matches the sheet's custom HyperLink style

# ---------------------------------------------------------------------
# Overview sheet: just the file name (via hyperlink) + status text
# ---------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")

$wsOverview.Hyperlinks.Delete()
$wsOverview.Hyperlinks.Add($wsOverview.Range("A2"), $mdUrl, [System.Reflection.Missing]::Value, [System.Reflection.Missing]::Value, $newGuidMd) | Out-Null
$wsOverview.Hyperlinks.Add($wsOverview.Range("A3"), $configUrl, [System.Reflection.Missing]::Value, [System.Reflection.Missing]::Value, ".localization-config") | Out-Null
$wsOverview.Range("A2").Font.Underline = $true
$wsOverview.Range("A2").Font.Color = $hlColor
$wsOverview.Range("A3").Font.Underline = $true
$wsOverview.Range("A3").Font.Color = $hlColor

$wsOverview.Range("B2").Value = $newStatus
$wsOverview.Range("C2").Value = $newStatus

# ---------------------------------------------------------------------
# Per-language detail sheets (zh-cn, de-de): the handoff attempt failed
# before a target file/date could be produced, so the handoff-file (C),
# handback-datetime (G) and handoff-reason (H) columns reset/collapse.
# ---------------------------------------------------------------------
foreach ($sheetName in @("zh-cn", "de-de")) {
    $ws = $wb.Worksheets.Item($sheetName)

    $ws.Hyperlinks.Delete()
    $ws.Hyperlinks.Add($ws.Range("A2"), $mdUrl, [System.Reflection.Missing]::Value, [System.Reflection.Missing]::Value, $newGuidMd) | Out-Null
    $ws.Hyperlinks.Add($ws.Range("A3"), $configUrl, [System.Reflection.Missing]::Value, [System.Reflection.Missing]::Value, ".localization-config") | Out-Null
    $ws.Range("A2").Font.Underline = $true
    $ws.Range("A2").Font.Color = $hlColor
    $ws.Range("A3").Font.Underline = $true
    $ws.Range("A3").Font.Color = $hlColor

    $ws.Range("B2").Value = $newStatus

    # The "Latest Handoff File" hyperlink/cell (C2) no longer applies.
    $ws.Range("C2").Clear()

    # Handoff/handback datetimes collapse to the empty-date sentinel.
    $ws.Range("D2").Value = $noDate
    $ws.Range("D3").Value = $noDate
    $ws.Range("G2").Value = $noDate
    $ws.Range("G3").Value = $noDate

    # Handoff reason becomes "Ignored" for both rows.
    $ws.Range("H2").Value = $ignored
    $ws.Range("H3").Value = $ignored
}
